$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.220.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.36%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.320.09'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.07%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.28%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.50'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.17%  '

$ws.Range('E7').Value = '  +0.24%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.574'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.67%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.318.03'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.06%  '

$ws.Range('E10').Value = '  -0.17%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.49'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.79%  '

$ws.Range('E12').Value = '  -0.49%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.333'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.06%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.736.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.94%  '

$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.238.10'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.13%  '

$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.48%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.39%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.325.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.56%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.41%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '313.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.16%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.70%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.61%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('E24').Value = '  +0.39%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.171'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.22%  '

$ws.Range('E26').Value = '  +0.29%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.51%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.01%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.44%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.91'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.15%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.42%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0726'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.98%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.89'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.64%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.36'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.381'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.06%  '

$ws.Range('E36').Value = '  +0.00%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.58%  '

$ws.Range('E38').Value = '  +0.30%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.09'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.67%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '315.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.03%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.06'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.72%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.72%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.86%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.13%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0938'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.93%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.569'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.44%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '18.84'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.12%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0492'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.10%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0220'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.18%  '

$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0213'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.97%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.17%  '
